$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9999
$ws.Range("I18").Value = 9999
$ws.Range("K18").Value = 9999
$ws.Range("M18").Value = -9715
$ws.Range("H33").Value = 628.3333
$ws.Range("I33").Value = 838.1111
$ws.Range("K33").Value = 838.1111
$ws.Range("M33").Value = -609.1111
$ws.Range("H53").Value = 5193.2383
$ws.Range("I53").Value = 257.92307
$ws.Range("K53").Value = 257.92307
$ws.Range("M53").Value = 379.07693
$ws.Range("H76").Value = 3477504.5
$ws.Range("I76").Value = 9262973
$ws.Range("K76").Value = 9262973
$ws.Range("M76").Value = -9262658
$ws.Range("H79").Value = 3477504.5
$ws.Range("I79").Value = 9262973
$ws.Range("K79").Value = 9262973
$ws.Range("M79").Value = -9261881
$ws.Range("H97").Value = 2679.8
$ws.Range("J97").Value = 2679.8
$ws.Range("L97").Value = 8039.400000000001
$ws.Range("N97").Value = -9031.400000000001
$ws.Range("H112").Value = 4501.3335
$ws.Range("J112").Value = 4967.4688
$ws.Range("L112").Value = 14902.4064
$ws.Range("N112").Value = -17118.4064
$ws.Range("H132").Value = 4547.278
$ws.Range("I132").Value = 4546.8237
$ws.Range("K132").Value = 13640.4711
$ws.Range("M132").Value = -11110.4711
$ws.Range("H137").Value = 53142.258
$ws.Range("I137").Value = 73199.24000000001
$ws.Range("J137").Value = 2999.8
$ws.Range("K137").Value = 219597.72
$ws.Range("L137").Value = 8999.400000000001
$ws.Range("M137").Value = -217047.72
$ws.Range("N137").Value = -14099.4
$ws.Range("H138").Value = 3180.258
$ws.Range("J138").Value = 3589.78
$ws.Range("L138").Value = 10769.34
$ws.Range("N138").Value = -21049.34

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 11530.597
$ws.Range("I32").Value = 7424.025
$ws.Range("K32").Value = 7424.025
$ws.Range("M32").Value = -7137.025
$ws.Range("H45").Value = 49222.09
$ws.Range("I45").Value = 78563.38
$ws.Range("K45").Value = 78563.38
$ws.Range("M45").Value = -78186.38
$ws.Range("H61").Value = 3153.4
$ws.Range("I61").Value = 3001.1875
$ws.Range("K61").Value = 3001.1875
$ws.Range("M61").Value = -2789.1875
$ws.Range("H96").Value = 67449.5
$ws.Range("J96").Value = 67449.5
$ws.Range("L96").Value = 67449.5
$ws.Range("N96").Value = -72941.5
$ws.Range("H97").Value = 8527.916999999999
$ws.Range("I97").Value = 20974.75
$ws.Range("K97").Value = 20974.75
$ws.Range("M97").Value = -20478.75
$ws.Range("H122").Value = 10586918
$ws.Range("I122").Value = 20205026
$ws.Range("K122").Value = 60615078
$ws.Range("M122").Value = -60612628
$ws.Range("H132").Value = 2899.8538
$ws.Range("I132").Value = 2276.6453
$ws.Range("K132").Value = 6829.9359
$ws.Range("M132").Value = -4299.9359
$ws.Range("H136").Value = 3153.4
$ws.Range("I136").Value = 3001.1875
$ws.Range("K136").Value = 9003.5625
$ws.Range("M136").Value = -6453.5625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H86").Value = 4666.967
$ws.Range("I86").Value = 5587.0454
$ws.Range("J86").Value = 2136.75
$ws.Range("K86").Value = 5587.0454
$ws.Range("L86").Value = 2136.75
$ws.Range("M86").Value = -4464.0454
$ws.Range("N86").Value = -4382.75
$ws.Range("H89").Value = 4666.967
$ws.Range("I89").Value = 5587.0454
$ws.Range("J89").Value = 2136.75
$ws.Range("K89").Value = 27935.227
$ws.Range("L89").Value = 10683.75
$ws.Range("M89").Value = -22319.227
$ws.Range("N89").Value = -21915.75
$ws.Range("H94").Value = 6483.6665
$ws.Range("I94").Value = 2203.625
$ws.Range("K94").Value = 2203.625
$ws.Range("M94").Value = -1752.625
$ws.Range("H107").Value = 2505.15
$ws.Range("I107").Value = 2236.1765
$ws.Range("J107").Value = 4029.3333
$ws.Range("K107").Value = 2236.1765
$ws.Range("L107").Value = 4029.3333
$ws.Range("M107").Value = -316.1765
$ws.Range("N107").Value = -7869.3333
$ws.Range("H134").Value = 3715.1082
$ws.Range("I134").Value = 1637
$ws.Range("K134").Value = 4911
$ws.Range("M134").Value = -2376

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3039.0469
$ws.Range("I58").Value = 4841.2964
$ws.Range("K58").Value = 4841.2964
$ws.Range("M58").Value = -4638.2964
$ws.Range("H86").Value = 11756.77
$ws.Range("I86").Value = 11289.2
$ws.Range("K86").Value = 11289.2
$ws.Range("M86").Value = -10166.2
$ws.Range("H89").Value = 11756.77
$ws.Range("I89").Value = 11289.2
$ws.Range("K89").Value = 56446
$ws.Range("M89").Value = -50830
$ws.Range("H94").Value = 1094.7142
$ws.Range("J94").Value = 1749.5
$ws.Range("L94").Value = 1749.5
$ws.Range("N94").Value = -2651.5
$ws.Range("H99").Value = 4082.2104
$ws.Range("I99").Value = 3825.7273
$ws.Range("J99").Value = 4434.875
$ws.Range("K99").Value = 3825.7273
$ws.Range("L99").Value = 4434.875
$ws.Range("M99").Value = -2327.7273
$ws.Range("N99").Value = -7430.875
$ws.Range("H109").Value = 19411.2
$ws.Range("J109").Value = 19889
$ws.Range("L109").Value = 19889
$ws.Range("N109").Value = -21969
$ws.Range("H126").Value = 4082.2104
$ws.Range("I126").Value = 3825.7273
$ws.Range("J126").Value = 4434.875
$ws.Range("K126").Value = 11477.1819
$ws.Range("L126").Value = 13304.625
$ws.Range("M126").Value = -9007.1819
$ws.Range("N126").Value = -18244.625
$ws.Range("H132").Value = 50958.973
$ws.Range("I132").Value = 35182.465
$ws.Range("J132").Value = 103547.336
$ws.Range("K132").Value = 105547.395
$ws.Range("L132").Value = 310642.008
$ws.Range("M132").Value = -103017.395
$ws.Range("N132").Value = -315702.008
$ws.Range("H134").Value = 26486.842
$ws.Range("I134").Value = 32824.414
$ws.Range("K134").Value = 98473.242
$ws.Range("M134").Value = -95938.242
$ws.Range("H136").Value = 3039.0469
$ws.Range("I136").Value = 4841.2964
$ws.Range("K136").Value = 14523.8892
$ws.Range("M136").Value = -11973.8892

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 68634.84
$ws.Range("I12").Value = 148363.67
$ws.Range("K12").Value = 445091.01
$ws.Range("M12").Value = -444918.01
$ws.Range("H17").Value = 148.23077
$ws.Range("I17").Value = 89.833336
$ws.Range("K17").Value = 269.500008
$ws.Range("M17").Value = -100.500008
$ws.Range("H22").Value = 1994
$ws.Range("J22").Value = 1991.5
$ws.Range("L22").Value = 5974.5
$ws.Range("N22").Value = -6312.5
$ws.Range("H25").Value = 510.31818
$ws.Range("I25").Value = 493.33334
$ws.Range("J25").Value = 546.7143
$ws.Range("K25").Value = 1480.00002
$ws.Range("L25").Value = 1640.1429
$ws.Range("M25").Value = -1311.00002
$ws.Range("N25").Value = -1978.1429
$ws.Range("H27").Value = 1994
$ws.Range("J27").Value = 1991.5
$ws.Range("L27").Value = 5974.5
$ws.Range("N27").Value = -6178.5
$ws.Range("H30").Value = 510.31818
$ws.Range("I30").Value = 493.33334
$ws.Range("J30").Value = 546.7143
$ws.Range("K30").Value = 1480.00002
$ws.Range("L30").Value = 1640.1429
$ws.Range("M30").Value = -1378.00002
$ws.Range("N30").Value = -1844.1429
$ws.Range("H39").Value = 3158
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 32258658
$ws.Range("I2").Value = 773.26666
$ws.Range("K2").Value = 773.26666
$ws.Range("M2").Value = -660.26666
$ws.Range("H12").Value = 4999.3335
$ws.Range("I12").Value = 4999
$ws.Range("J12").Value = 4999.5
$ws.Range("K12").Value = 4999
$ws.Range("L12").Value = 4999.5
$ws.Range("M12").Value = -4859
$ws.Range("N12").Value = -5279.5
$ws.Range("H14").Value = 33667332
$ws.Range("J14").Value = 100000000
$ws.Range("L14").Value = 100000000
$ws.Range("N14").Value = -100000336
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -3471
$ws.Range("N22").Value = -5058
$ws.Range("H39").Value = 59994
$ws.Range("J39").Value = 59994
$ws.Range("L39").Value = 59994
$ws.Range("N39").Value = -61058
$ws.Range("H70").Value = 7090.909
$ws.Range("I70").Value = 6427.4287
$ws.Range("K70").Value = 6427.4287
$ws.Range("M70").Value = -6157.4287
$ws.Range("H73").Value = 7090.909
$ws.Range("I73").Value = 6427.4287
$ws.Range("K73").Value = 6427.4287
$ws.Range("M73").Value = -5491.4287
$ws.Range("H126").Value = 4062.516
$ws.Range("I126").Value = 2549.4443
$ws.Range("J126").Value = 4681.5
$ws.Range("K126").Value = 7648.3329
$ws.Range("L126").Value = 14044.5
$ws.Range("M126").Value = -5178.3329
$ws.Range("N126").Value = -18984.5
$ws.Range("H132").Value = 2465.182
$ws.Range("I132").Value = 2384.075
$ws.Range("K132").Value = 7152.224999999999
$ws.Range("M132").Value = -4622.224999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6067.467
$ws.Range("I7").Value = 3520.6
$ws.Range("J7").Value = 11161.2
$ws.Range("K7").Value = 3520.6
$ws.Range("L7").Value = 11161.2
$ws.Range("M7").Value = -3408.6
$ws.Range("N7").Value = -11385.2
$ws.Range("H22").Value = 69394.766
$ws.Range("I22").Value = 127804.71
$ws.Range("K22").Value = 127804.71
$ws.Range("M22").Value = -127509.71
$ws.Range("H27").Value = 69394.766
$ws.Range("I27").Value = 127804.71
$ws.Range("K27").Value = 127804.71
$ws.Range("M27").Value = -127697.71
$ws.Range("H55").Value = 2479.3572
$ws.Range("J55").Value = 1859.4
$ws.Range("L55").Value = 1859.4
$ws.Range("N55").Value = -2205.4
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H126").Value = 6067.467
$ws.Range("I126").Value = 3520.6
$ws.Range("J126").Value = 11161.2
$ws.Range("K126").Value = 10561.8
$ws.Range("L126").Value = 33483.60000000001
$ws.Range("M126").Value = -8091.799999999999
$ws.Range("N126").Value = -38423.60000000001
$ws.Range("H132").Value = 10084.85
$ws.Range("I132").Value = 11120.467
$ws.Range("J132").Value = 6978
$ws.Range("K132").Value = 33361.401
$ws.Range("L132").Value = 20934
$ws.Range("M132").Value = -30831.401
$ws.Range("N132").Value = -25994

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 190888.8
$ws.Range("J4").Value = 7222
$ws.Range("L4").Value = 7222
$ws.Range("N4").Value = -7448
$ws.Range("H95").Value = 51349
$ws.Range("J95").Value = 51349
$ws.Range("L95").Value = 51349
$ws.Range("N95").Value = -56841
$ws.Range("H109").Value = 67000
$ws.Range("J109").Value = 67000
$ws.Range("L109").Value = 67000
$ws.Range("N109").Value = -69774
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 2990.8
$ws.Range("I126").Value = 3086.8572
$ws.Range("K126").Value = 9260.571599999999
$ws.Range("M126").Value = -6790.571599999999
$ws.Range("H132").Value = 230792.48
$ws.Range("I132").Value = 4049.1365
$ws.Range("K132").Value = 12147.4095
$ws.Range("M132").Value = -9617.4095
